$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "34.408.88"
Set-TextValue "E2" "  -0.42%  "
Set-TextValue "D3" "1.800.70"
Set-TextValue "E3" "  -0.90%  "
Set-TextValue "E4" "  -0.37%  "
Set-TextValue "D5" "225.43"
Set-TextValue "E5" "  -1.38%  "
Set-TextValue "E6" "  +3.58%  "
Set-TextValue "E7" "  -0.34%  "
Set-TextValue "D8" "36.24"
Set-TextValue "E8" "  +3.43%  "
Set-TextValue "E9" "  -2.97%  "
Set-TextValue "E10" "  -2.80%  "
Set-TextValue "D11" "0.0964"
Set-TextValue "E11" "  +1.11%  "
Set-TextValue "D12" "2.059.16"
Set-TextValue "E12" "  -0.96%  "
Set-TextValue "D13" "11.22"
Set-TextValue "E13" "  -1.25%  "
Set-TextValue "D14" "1.786.15"
Set-TextValue "E14" "  -1.57%  "
Set-TextValue "E15" "  -2.41%  "
Set-TextValue "D16" "34.378.90"
Set-TextValue "E16" "  -0.50%  "
Set-TextValue "D17" "4.41"
Set-TextValue "E17" "  +1.08%  "
Set-TextValue "D18" "68.75"
Set-TextValue "E18" "  -0.95%  "
Set-TextValue "D19" "245.64"
Set-TextValue "E19" "  -0.11%  "
Set-TextValue "E20" "  -3.49%  "
Set-TextValue "D21" "11.36"
Set-TextValue "E21" "  -2.10%  "
Set-TextValue "E22" "  -0.40%  "
Set-TextValue "D23" "4.07"
Set-TextValue "E23" "  -2.85%  "
Set-TextValue "D24" "2.19"
Set-TextValue "E24" "  +4.21%  "
Set-TextValue "D25" "170.67"
Set-TextValue "E25" "  -0.77%  "
Set-TextValue "D26" "7.89"
Set-TextValue "E26" "  +4.44%  "
Set-TextValue "D27" "17.37"
Set-TextValue "E27" "  +3.22%  "
Set-TextValue "D28" "0.120"
Set-TextValue "E28" "  +1.63%  "
Set-TextValue "E29" "  -0.39%  "
Set-TextValue "E30" "  -1.92%  "
Set-TextValue "D31" "3.78"
Set-TextValue "E31" "  -1.49%  "
Set-TextValue "D32" "3.89"
Set-TextValue "E32" "  -2.44%  "
Set-TextValue "D33" "0.0513"
Set-TextValue "E33" "  -3.22%  "
Set-TextValue "E34" "  -4.31%  "
Set-TextValue "D35" "1.360.17"
Set-TextValue "E35" "  -3.06%  "
Set-TextValue "D36" "0.646"
Set-TextValue "E36" "  -5.31%  "
Set-TextValue "E37" "  -1.35%  "
Set-TextValue "D38" "2.35"
Set-TextValue "E38" "  -8.71%  "
Set-TextValue "D39" "0.0186"
Set-TextValue "E39" "  -2.62%  "
Set-TextValue "D40" "2.42"
Set-TextValue "E40" "  +0.40%  "
Set-TextValue "E41" "  -2.44%  "
Set-TextValue "D42" "81.03"
Set-TextValue "E42" "  -2.74%  "
Set-TextValue "D43" "0.934"
Set-TextValue "E43" "  -2.19%  "
Set-TextValue "E44" "  +4.95%  "
Set-TextValue "D45" "13.16"
Set-TextValue "E45" "  -4.72%  "
Set-TextValue "E46" "  -2.76%  "
Set-TextValue "D47" "1.962.54"
Set-TextValue "E47" "  -0.85%  "
Set-TextValue "D48" "5.78"
Set-TextValue "E48" "  -4.54%  "
Set-TextValue "E49" "  -0.38%  "
Set-TextValue "D50" "101.56"
Set-TextValue "E50" "  -3.85%  "
Set-TextValue "E51" "  -7.68%  "
